$wb = $excel.ActiveWorkbook

# --- 1. Insert a new worksheet "2022-Q3" right after "总计" (before "2022-Q2") ---
$wsTotal = $wb.Worksheets.Item("总计")
$wsOldQ2 = $wb.Worksheets.Item("2022-Q2")
$wsNew = $wb.Worksheets.Add($wsOldQ2)
$wsNew.Name = "2022-Q3"

function Set-HeaderStyle($rng) {
    $rng.Font.Bold = $true
    $rng.HorizontalAlignment = -4108
    $rng.VerticalAlignment = -4160
    $rng.Borders.LineStyle = 1
}

# --- 2. Populate the new "2022-Q3" sheet with fund-holding data ---
$wsNew.Range("B1").Value = "基金代码"
$wsNew.Range("C1").Value = "基金名称"
$wsNew.Range("D1").Value = "基金规模"
$wsNew.Range("E1").Value = "股票总仓位"
$wsNew.Range("F1").Value = "仓位占比"
$wsNew.Range("G1").Value = "持有市值(亿元)"
$wsNew.Range("H1").Value = "仓位排名"
Set-HeaderStyle $wsNew.Range("B1:H1")

$wsNew.Range("A2").Value = 0
$wsNew.Range("B2").Value = "008545"
$wsNew.Range("C2").Value = "泓德丰润三年持有期混合"
$wsNew.Range("D2").Value = "74.71"
$wsNew.Range("E2").Value = "88.42"
$wsNew.Range("F2").Value = "4.53"
$wsNew.Range("G2").Value = "3.3844"
$wsNew.Range("H2").Value = 6
Set-HeaderStyle $wsNew.Range("A2")

$wsNew.Range("A3").Value = 1
$wsNew.Range("B3").Value = "005395"
$wsNew.Range("C3").Value = "泓德臻远回报灵活配置混合"
$wsNew.Range("D3").Value = "29.76"
$wsNew.Range("E3").Value = "93.35"
$wsNew.Range("F3").Value = "5.48"
$wsNew.Range("G3").Value = "1.6308"
$wsNew.Range("H3").Value = 6
Set-HeaderStyle $wsNew.Range("A3")

$wsNew.Range("A4").Value = 2
$wsNew.Range("B4").Value = "010864"
$wsNew.Range("C4").Value = "泓德卓远混合A"
$wsNew.Range("D4").Value = "22.84"
$wsNew.Range("E4").Value = "92.87"
$wsNew.Range("F4").Value = "3.71"
$wsNew.Range("G4").Value = "0.8474"
$wsNew.Range("H4").Value = 7
Set-HeaderStyle $wsNew.Range("A4")

$wsNew.Range("A5").Value = 3
$wsNew.Range("B5").Value = "010865"
$wsNew.Range("C5").Value = "泓德卓远混合C"
$wsNew.Range("D5").Value = "10.32"
$wsNew.Range("E5").Value = "92.87"
$wsNew.Range("F5").Value = "3.71"
$wsNew.Range("G5").Value = "0.3829"
$wsNew.Range("H5").Value = 7
Set-HeaderStyle $wsNew.Range("A5")

# --- 3. Insert a new row into "总计" for the 2022-Q3 summary, shifting the rest down ---
$wsTotal.Range("A3:D7").Insert()

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 4
$wsTotal.Range("D2").Value = 6.25
Set-HeaderStyle $wsTotal.Range("A2")

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q2"
$wsTotal.Range("C3").Value = 3
$wsTotal.Range("D3").Value = 3.64
Set-HeaderStyle $wsTotal.Range("A3")

$wsTotal.Range("A4").Value = 2
$wsTotal.Range("B4").Value = "2021-Q3"
$wsTotal.Range("C4").Value = 13
$wsTotal.Range("D4").Value = 11.65
Set-HeaderStyle $wsTotal.Range("A4")

$wsTotal.Range("A5").Value = 3
$wsTotal.Range("B5").Value = "2021-Q2"
$wsTotal.Range("C5").Value = 5
$wsTotal.Range("D5").Value = 5.32
Set-HeaderStyle $wsTotal.Range("A5")

$wsTotal.Range("A6").Value = 4
$wsTotal.Range("B6").Value = "2021-Q1"
$wsTotal.Range("C6").Value = 11
$wsTotal.Range("D6").Value = 6.54
Set-HeaderStyle $wsTotal.Range("A6")

$wsTotal.Range("A7").Value = 5
$wsTotal.Range("B7").Value = "2020-Q4"
$wsTotal.Range("C7").Value = 12
$wsTotal.Range("D7").Value = 10.45
Set-HeaderStyle $wsTotal.Range("A7")
